# "Committed validations in personalise page"
#
# The sheet ("Status") originally held a 2-row data table (rows 2-3) plus a
# long tail of 33 completely empty, but still "dimensioned", rows (4-36)
# left over from a bigger sheet that used to live here. The edit:
#   1. drops all of that dead trailing range so the sheet's used range
#      shrinks back down to the real data (A1:G3),
#   2. swaps which "Mid" value sits in row 2 vs row 3 (and updates the text
#      of one of the two Mid values), and
#   3. moves the active selection from the old out-of-range cell (J4) onto
#      the new last row of data (A3), clearing the stale topLeftCell scroll
#      position at the same time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the empty trailing rows (4-36) entirely so the sheet's used
#    range / dimension collapses back to the real A1:G3 data block.
$ws.Rows("4:36").Delete()

# 2. Row 2's Mid becomes "20150914005" (previously row 3's "20150914007"
#    text, itself corrected), row 3's Mid becomes "20150914006" (what used
#    to be row 2's value). Only column A changes - the rest of each row
#    (Used/SAQ Type/SAQ Complete?/Attested?/Scenario/Username) is untouched.
$ws.Range("A2").Value = "20150914005"
$ws.Range("A3").Value = "20150914006"

# 3. Select A3 (the new bottom data row) instead of the old J4 selection.
$ws.Range("A3").Select()
